# Applies the "changes to event table 2.0" commit:
#  - wraps "EVENT  TABLE" in grammar proofErr marks
#  - adds spellStart/spellEnd proofErr marks around several misspelled /
#    foreign words, splitting their runs accordingly
#  - fixes "While opponent is in the air Jump+Grapple" -> "While enemy is in
#    the air Jump+Grapple" (plus spell marks)
#  - adds a comma after "Hold Crouch(Key Board CTRL 360 Pad LB)"
#  - moves the "_GoBack" bookmark from the empty trailing "Game" cell to
#    wrap the "While Dashing..." cell text

$d = $word.ActiveDocument
$W_NS = 'http://schemas.openxmlformats.org/wordprocessingml/2006/main'

function Replace-WholeParagraph {
    param(
        [string]$OldText,
        [string]$InnerXml
    )

    $r = $d.Content
    $found = $r.Find.Execute($OldText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Output "NOT FOUND: $OldText"
        return $null
    }
    $xml = '<w:p xmlns:w="' + $W_NS + '">' + $InnerXml + '</w:p>'
    $r.InsertXML($xml)
    Write-Output "Replaced: $OldText"
    return $r
}

# 1. "EVENT  TABLE" heading: wrap the bold/underlined run in gramStart/gramEnd
Replace-WholeParagraph `
    "EVENT  TABLE:  List all the events that you anticipate your software will do and fill in the rest of the columns" `
    (
        '<w:proofErr w:type="gramStart"/>' +
        '<w:r><w:rPr><w:b/><w:u w:val="single"/></w:rPr><w:t>EVENT  TABLE</w:t></w:r>' +
        '<w:proofErr w:type="gramEnd"/>' +
        '<w:r><w:t>:  List all the events that you anticipate your software will do and fill in the rest of the columns</w:t></w:r>'
    )

# 2. "Keyboard (Space Bar) Cotroller (360 Pad RT)" -> spellcheck mark "Cotroller"
Replace-WholeParagraph `
    "Keyboard (Space Bar) Cotroller (360 Pad RT)" `
    (
        '<w:r><w:t xml:space="preserve">Keyboard (Space Bar) </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/><w:r><w:t>Cotroller</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t xml:space="preserve"> (360 Pad RT)</w:t></w:r>'
    )

# 3. "...for an area of effect (AoE)" -> spellcheck mark "AoE"
Replace-WholeParagraph `
    "Character air dive and punches the ground for an area of effect (AoE)" `
    (
        '<w:r><w:t>Character air dive and punches the ground for an area of effect (</w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/><w:r><w:t>AoE</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t>)</w:t></w:r>'
    )

# 4. "Character does Meia Lua de Compasso Dupla " -> spellcheck marks on each Portuguese word
Replace-WholeParagraph `
    "Character does Meia Lua de Compasso Dupla " `
    (
        '<w:r><w:t xml:space="preserve">Character does </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/><w:r><w:t>Meia</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/><w:r><w:t>Lua</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t xml:space="preserve"> de </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/><w:r><w:t>Compasso</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/><w:r><w:t>Dupla</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t xml:space="preserve"> </w:t></w:r>'
    )

# 5. "Character does bencao" -> spellcheck mark "bencao"
Replace-WholeParagraph `
    "Character does bencao" `
    (
        '<w:r><w:t xml:space="preserve">Character does </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/><w:r><w:t>bencao</w:t></w:r><w:proofErr w:type="spellEnd"/>'
    )

# 6. "Light Kick (Meia-lua de Frente) to give distance..." -> spellcheck marks on "Meia-lua" / "Frente"
Replace-WholeParagraph `
    "Light Kick (Meia-lua de Frente) to give distance to the enemy form the user (up direction)" `
    (
        '<w:r><w:t>Light Kick (</w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/><w:r><w:t>Meia-lua</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t xml:space="preserve"> de </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/><w:r><w:t>Frente</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t>) to give distance to the enemy form the user (up direction)</w:t></w:r>'
    )

# 7. "Before the Corkscew animation finish the character grabs..." -> spellcheck mark "Corkscew"
Replace-WholeParagraph `
    "Before the Corkscew animation finish the character grabs the enemy while using the rotation on an enemy body as a weapon until stamina has been drain can move with this move " `
    (
        '<w:r><w:t xml:space="preserve">Before the </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/><w:r><w:t>Corkscew</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t xml:space="preserve"> animation finish the character </w:t></w:r>' +
        '<w:r><w:t xml:space="preserve">grabs the enemy </w:t></w:r>' +
        '<w:r><w:t xml:space="preserve">while using the rotation </w:t></w:r>' +
        '<w:r><w:t>on an enemy body as a weapon until stamina has been drain</w:t></w:r>' +
        '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
        '<w:r><w:t>can move with this move</w:t></w:r>' +
        '<w:r><w:t xml:space="preserve"> </w:t></w:r>'
    )

# 8. "While opponent is in the air Jump+Grapple" -> "While enemy is in the air Jump+Grapple" (+ spell marks)
Replace-WholeParagraph `
    "While opponent is in the air Jump+Grapple" `
    (
        '<w:r><w:t xml:space="preserve">While </w:t></w:r>' +
        '<w:r><w:t>enemy</w:t></w:r>' +
        '<w:r><w:t xml:space="preserve"> is in the air </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:t>Jump+</w:t></w:r>' +
        '<w:r><w:t>G</w:t></w:r>' +
        '<w:r><w:t>rapple</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>'
    )

# 9. "Hold Crouch(Key Board CTRL 360 Pad LB) " (Hammer Time row) -> add comma after the closing paren
Replace-WholeParagraph `
    "Hold Crouch(Key Board CTRL 360 Pad LB)  Move left and right twice, Release Crouch(Key Board CTRL 360 Pad LB)  + Heavy Punch( Keyboard I, 360 Pad  Y)" `
    (
        '<w:r><w:t>Hold Crouch(Key Board CTRL 360 Pad LB)</w:t></w:r>' +
        '<w:r><w:t>,</w:t></w:r>' +
        '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
        '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
        '<w:r><w:t>Move left and right twice,</w:t></w:r>' +
        '<w:r><w:t xml:space="preserve"> Release</w:t></w:r>' +
        '<w:r><w:t xml:space="preserve"> Crouch(Key Board CTRL 360 Pad LB) </w:t></w:r>' +
        '<w:r><w:t xml:space="preserve"> +</w:t></w:r>' +
        '<w:r><w:t xml:space="preserve"> Heavy Punch</w:t></w:r>' +
        '<w:r><w:t>( Keyboard I, 360 Pad  Y)</w:t></w:r>'
    )

# 10. Move the "_GoBack" bookmark from the trailing "Game" cell onto the
#     "While Dashing..." cell text.
Replace-WholeParagraph `
    "While Dashing Hold Crouch(Key Board CTRL 360 Pad LB) + Direction you were Dashing" `
    (
        '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
        '<w:r><w:t xml:space="preserve">While Dashing </w:t></w:r>' +
        '<w:r><w:t>Hold Cr</w:t></w:r>' +
        '<w:r><w:t>ouch(Key Board CTRL 360 Pad LB) + Direction you were Dashing</w:t></w:r>' +
        '<w:bookmarkEnd w:id="0"/>'
    )

# Now find the "Game" cell that immediately follows (still has the old
# bookmark) and strip it back down to a plain run.
$r = $d.Content
$found1 = $r.Find.Execute("Dashes and Roll into a giant disco ball doing damage those in the way", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $r.Collapse(0)
    $found2 = $r.Find.Execute("Game", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found2) {
        $xml = '<w:p xmlns:w="' + $W_NS + '"><w:r><w:t>Game</w:t></w:r></w:p>'
        $r.InsertXML($xml)
        Write-Output "Replaced: trailing Game/_GoBack cell"
    } else {
        Write-Output "NOT FOUND: trailing Game cell"
    }
} else {
    Write-Output "NOT FOUND: Dashes and Roll... cell"
}

# 11. "Select Strenght" -> spellcheck mark "Strenght"
Replace-WholeParagraph `
    "Select Strenght" `
    (
        '<w:r><w:t xml:space="preserve">Select </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/><w:r><w:t>Strenght</w:t></w:r><w:proofErr w:type="spellEnd"/>'
    )
